# Upload new version with timestamp
# Adds three new shortage-report rows (COLOVERIN, NEVILOB, VASTAREL) into the
# alphabetically sorted product list, renumbers the "م" sequence column,
# updates the running total and refreshes the printed timestamp.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Insert three new data rows at the correct alphabetical positions.
#    Each insert copies the row immediately above it (so number formats,
#    fonts, fills and borders carry over) and shifts everything below down.
# ---------------------------------------------------------------------------

# COLOVERIN A 30 TABLETS  -> goes right after COLOVATIL (row 11), before CYANOHEPTAN (row 12)
$ws.Rows("11").Copy()
$ws.Rows("12").Insert()

# NEVILOB 2.5 MG 14 TAB.  -> goes right after MOTINORM (row 15), before OXALEPTAL (now row 16)
$ws.Rows("15").Copy()
$ws.Rows("16").Insert()

# VASTAREL MR 35MG 30 F.C.TAB. -> goes right after TREFLUCAN (now row 23), before سرنجات (now row 24)
$ws.Rows("23").Copy()
$ws.Rows("24").Insert()

# ---------------------------------------------------------------------------
# 2) Re-create the merged cells on the three freshly inserted rows (Insert()
#    does not copy merge information from the source row).
# ---------------------------------------------------------------------------
foreach ($r in 12, 16, 24) {
    $ws.Range("A$r`:B$r").Merge() | Out-Null
    $ws.Range("C$r`:G$r").Merge() | Out-Null
    $ws.Range("H$r`:K$r").Merge() | Out-Null
    $ws.Range("L$r`:M$r").Merge() | Out-Null
    $ws.Range("N$r`:O$r").Merge() | Out-Null
}

# ---------------------------------------------------------------------------
# 3) Fill in the cell values for the new rows.
# ---------------------------------------------------------------------------
$ws.Range("C12").Value = "COLOVERIN A 30 TABLETS"
$ws.Range("H12").Value = "1:2"
$ws.Range("L12").Value = "1"
$ws.Range("N12").Value = "99.00"
$ws.Range("P12").Value = "32.6700"
$ws.Range("Q12").Value = "0:1"

$ws.Range("C16").Value = "NEVILOB 2.5 MG 14 TAB."
$ws.Range("H16").Value = "2:0"
$ws.Range("L16").Value = "1"
$ws.Range("N16").Value = "46.00"
$ws.Range("P16").Value = "23.0000"
$ws.Range("Q16").Value = "0:1"

$ws.Range("C24").Value = "VASTAREL MR 35MG 30 F.C.TAB."
$ws.Range("H24").Value = "1:0"
$ws.Range("L24").Value = "1"
$ws.Range("N24").Value = "175.00"
$ws.Range("P24").Value = "57.7500"
$ws.Range("Q24").Value = "0:1"

# ---------------------------------------------------------------------------
# 4) Renumber the "م" (sequence) column for all 20 data rows (7-26).
# ---------------------------------------------------------------------------
$seq = 1
for ($r = 7; $r -le 26; $r++) {
    $ws.Range("A$r").Value = $seq
    $seq++
}

# ---------------------------------------------------------------------------
# 5) Row heights: keep the same look Excel produced when the three rows were
#    inserted (matches the committed workbook's auto-fit result).
# ---------------------------------------------------------------------------
$heights = @{
    7=25.5; 8=24.75; 9=25.5; 10=24.75; 11=25.5; 12=25.5; 13=24.75; 14=25.5;
    15=24.75; 16=25.5; 17=25.5; 18=24.75; 19=25.5; 20=24.75; 21=25.5; 22=25.5;
    23=24.75; 24=25.5; 25=24.75; 26=25.5
}
foreach ($r in $heights.Keys) {
    $ws.Rows($r).RowHeight = $heights[$r]
}

# ---------------------------------------------------------------------------
# 6) Update the running total (sum of the "سعر البيع" column) and footer row.
# ---------------------------------------------------------------------------
$ws.Range("P27").Value = 832.54
$ws.Rows("27").RowHeight = 25.5

$ws.Range("A28").Value = "Wednesday, 27 August, 2025 11:25 AM"
$ws.Rows("28").RowHeight = 16.5
